$wb = $excel.ActiveWorkbook

# Markets that get the new MX-BBX / MX-DPBX accessory rows inserted just
# above the "Wg" section-break row on their "Panel Accessories" sheet.
$names = @("Slovakia", "Italy", "Netherlands", "Austria", "Denmark")

foreach ($name in $names) {
    $ws = $wb.Worksheets.Item($name)

    # Locate the "Wg" row (start of the next block) - the two new rows are
    # inserted immediately above it, pushing it (and everything below) down.
    $found = $ws.Columns.Item(1).Find("Wg")
    $row = $found.Row

    $insertRange = "A" + $row + ":A" + ($row + 1)
    $ws.Range($insertRange).EntireRow.Insert() | Out-Null

    $cellBBX = "A" + $row
    $cellDPBX = "A" + ($row + 1)
    $ws.Range($cellBBX).Value = "MX-BBX"
    $ws.Range($cellDPBX).Value = "MX-DPBX"

    # Pick up the formatting (borders etc.) of the row directly below, which
    # is the very row ("Wg") that just got displaced - keeps the same style
    # used throughout the accessories list instead of a blank/unstyled cell.
    $formatSource = "A" + ($row + 2)
    $ws.Range($formatSource).Copy() | Out-Null
    $ws.Range($insertRange).PasteSpecial(-4122) | Out-Null

    # Reflect the manual selection left behind on the two freshly added rows.
    $ws.Range($insertRange).Select() | Out-Null
}

# A couple of other tabs were simply clicked through while reviewing the
# above edits - their selections moved too, with no data change.
$wsSpain = $wb.Worksheets.Item("Spain")
$wsSpain.Range("A13:A14").Select() | Out-Null

$wsTurkey = $wb.Worksheets.Item("Turkey")
$wsTurkey.Range("A8:A15").Select() | Out-Null

# Denmark ends up the active tab.
$wsDenmark = $wb.Worksheets.Item("Denmark")
$wsDenmark.Activate() | Out-Null
$wsDenmark.Range("A14:A15").Select() | Out-Null
